$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Summary" (2nd worksheet)
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item(2)

$wsSummary.Range("B2").Value = 402.79
$wsSummary.Range("E2").Value = 9597.21
$wsSummary.Range("F2").Value = 1275.67
# G2 becomes a new, empty (default-styled) touched cell -> dimension grows to G5
$wsSummary.Range("G2").Borders.LineStyle = -4142

$wsSummary.Range("A5").Value = 0.23
$wsSummary.Range("B5").Value = 0.23

$wsSummary.Activate()
$wsSummary.Range("C5").Select()

# ---------------------------------------------------------------------------
# Sheet "Repayment schedule" (3rd worksheet)
# ---------------------------------------------------------------------------
$wsSchedule = $wb.Worksheets.Item(3)

# The stray empty P2 cell moves to O2 (O1 already has a header, P column is
# used by data rows further down so the sheet dimension stays A1:P15)
$wsSchedule.Range("P2").Clear()
$wsSchedule.Range("O2").Borders.LineStyle = -4142

$wsSchedule.Range("L3").Value = 476.76
$wsSchedule.Range("N3").Value = 476.76
$wsSchedule.Range("P3").Value = 410.96

$wsSchedule.Range("J5").Value = 0.23
$wsSchedule.Range("K5").Value = 887.95
$wsSchedule.Range("L5").Value = 23.24

$wsSchedule.Activate()
$wsSchedule.Range("D7").Select()

# ---------------------------------------------------------------------------
# Sheet "Transactions" (4th worksheet)
# ---------------------------------------------------------------------------
$wsTrans = $wb.Worksheets.Item(4)

$wsTrans.Range("A2").Value = 3578

$wsTrans.Range("A3").Value = 3577
$wsTrans.Range("E3").Value = 23.24
$wsTrans.Range("I3").Value = 0.23

# Rows 4 & 5 swap their Disbursement/Repayment data, including number formats.
# Re-point the formatting first (copy format only) before writing new values.
$wsTrans.Range("F4").Copy()
$wsTrans.Range("E4").PasteSpecial(-4122)
$wsTrans.Range("F4").Copy()
$wsTrans.Range("J4").PasteSpecial(-4122)
$wsTrans.Range("E7").Copy()
$wsTrans.Range("E5").PasteSpecial(-4122)
$wsTrans.Range("J7").Copy()
$wsTrans.Range("J5").PasteSpecial(-4122)

$wsTrans.Range("A4").Value = 3575
$wsTrans.Range("D4").Value = "Repayment"
$wsTrans.Range("E4").Value = 500
$wsTrans.Range("F4").Value = 402.79
$wsTrans.Range("G4").Value = 96.98
$wsTrans.Range("I4").Value = 0.23
$wsTrans.Range("J4").Value = 0

$wsTrans.Range("A5").Value = 3574
$wsTrans.Range("D5").Value = "Disbursement"
$wsTrans.Range("E5").Value = 5000
$wsTrans.Range("F5").Value = 0
$wsTrans.Range("G5").Value = 0
$wsTrans.Range("J5").Value = 10000

$wsTrans.Range("A6").Value = 3576

$wsTrans.Range("A7").Value = 3572

$wsTrans.Activate()
$wsTrans.Range("D7").Select()
